# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (cloned from "2021-Q4" so headers /
#    fund code / fund name / styles already match) positioned right before
#    the "总计" (totals) sheet, then refresh its numeric columns.
# 2. Refresh the "总计" sheet: shift every existing row down by one and
#    write the new 2022-Q1 summary row at the top (row 2).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the "2022-Q1" sheet from the "2021-Q4" template.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalBeforeCopy = $wb.Worksheets.Item("总计")
$template.Copy($totalBeforeCopy)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# NOTE: the worksheet-copy operation can leave earlier sheet-object
# references pointing at the wrong tab, so re-resolve "总计" by name now
# that the sheet collection has settled.
$total = $wb.Worksheets.Item("总计")

# Columns D:G look numeric but are stored as plain text in this workbook
# (note the significant trailing zero in "0.60"), so force text formatting
# before assigning, then drop back to the default style so no stray
# number-format style lingers on the cell.
$q1.Range("D2:G2").NumberFormat = "@"
$q1.Range("D2").Value = "0.60"
$q1.Range("E2").Value = "87.31"
$q1.Range("F2").Value = "8.72"
$q1.Range("G2").Value = "0.0523"
$q1.Range("D2:G2").Style = "Normal"

$q1.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 2. Refresh the "总计" sheet with the new quarter on top.
# ---------------------------------------------------------------------
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 1
$total.Range("D7").Value = 0.03

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.03

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.04

$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.05

$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.05

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.05

# New trailing row: give A7 the same index-column styling as A2:A6 by
# copying A6's format onto it (xlPasteFormats = -4122) after the value is
# set, so no stray style index gets introduced.
$total.Range("A7").Value = 5
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
